# The site footer that the old Jekyll build used to append after the
# "Requisitos" list - "Ver no Jupiter Salvar em pdf Salvar em docx" and
# the "(c) 2020 ... Creative Commons Attribution" copyright line - is
# gone from the freshly-built page, and with it one of the two blank
# paragraphs that used to sandwich it (the other blank paragraph, right
# before the trailing page-break paragraph, survives as-is).

$d = $word.ActiveDocument

$verPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Ver no Jupiter Salvar em pdf Salvar em docx*") {
        $verPara = $p
    }
}

if ($verPara -ne $null) {
    $leadBlankPara = $verPara.Previous()
    $copyrightPara = $verPara.Next()

    $deleteRange = $d.Range($leadBlankPara.Range.Start, $copyrightPara.Range.End)
    $deleteRange.Delete()
}
